$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1108.6757
$ws.Range("I17").Value = 1020.75
$ws.Range("J17").Value = 1119.3334
$ws.Range("K17").Value = 3062.25
$ws.Range("L17").Value = 3358.0002
$ws.Range("M17").Value = -2894.25
$ws.Range("N17").Value = -3694.0002

# Hunk 1: sheet ALC, row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 723.6667
$ws.Range("I33").Value = 310.25
$ws.Range("J33").Value = 1550.5
$ws.Range("K33").Value = 310.25
$ws.Range("L33").Value = 1550.5
$ws.Range("M33").Value = -81.25
$ws.Range("N33").Value = -2008.5

# Hunk 2: sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4710.357
$ws.Range("I40").Value = 4087.25
$ws.Range("J40").Value = 5541.1665
$ws.Range("K40").Value = 4087.25
$ws.Range("L40").Value = 5541.1665
$ws.Range("M40").Value = -3912.25
$ws.Range("N40").Value = -5891.1665

# Hunk 3: sheet ALC, row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 762.7143
$ws.Range("I41").Value = 355.66666
$ws.Range("J41").Value = 1068
$ws.Range("K41").Value = 355.66666
$ws.Range("L41").Value = 1068
$ws.Range("M41").Value = 84.33334000000002
$ws.Range("N41").Value = -1948

# Hunk 4: sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 14477.777
$ws.Range("I137").Value = 16573.732
$ws.Range("J137").Value = 3998
$ws.Range("K137").Value = 49721.196
$ws.Range("L137").Value = 11994
$ws.Range("M137").Value = -47171.196
$ws.Range("N137").Value = -17094

# Hunk 5: sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4704.5
$ws.Range("I138").Value = 2979.6667
$ws.Range("J138").Value = 5102.5386
$ws.Range("K138").Value = 8939.000100000001
$ws.Range("L138").Value = 15307.6158
$ws.Range("M138").Value = -3799.000100000001
$ws.Range("N138").Value = -25587.6158

# Hunk 6: sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6383.349
$ws.Range("I32").Value = 6272.1475
$ws.Range("K32").Value = 6272.1475
$ws.Range("M32").Value = -5985.1475

# Hunk 7: sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4784.2705
$ws.Range("I61").Value = 4374.769
$ws.Range("J61").Value = 5752.1816
$ws.Range("K61").Value = 4374.769
$ws.Range("L61").Value = 5752.1816
$ws.Range("M61").Value = -4162.769
$ws.Range("N61").Value = -6176.1816

# Hunk 8: sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7441.4
$ws.Range("I132").Value = 7836.206
$ws.Range("K132").Value = 23508.618
$ws.Range("M132").Value = -20978.618

# Hunk 9: sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4784.2705
$ws.Range("I136").Value = 4374.769
$ws.Range("J136").Value = 5752.1816
$ws.Range("K136").Value = 13124.307
$ws.Range("L136").Value = 17256.5448
$ws.Range("M136").Value = -10574.307
$ws.Range("N136").Value = -22356.5448

# Hunk 10: sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2886.658
$ws.Range("I105").Value = 2227.5667
$ws.Range("J105").Value = 5358.25
$ws.Range("K105").Value = 2227.5667
$ws.Range("L105").Value = 5358.25
$ws.Range("M105").Value = -480.5666999999999
$ws.Range("N105").Value = -8852.25

# Hunk 11: sheet BSM, row 139
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 73328.836
$ws.Range("J139").Value = 79994.60000000001
$ws.Range("L139").Value = 79994.60000000001
$ws.Range("N139").Value = -90274.60000000001

# Hunk 12: sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2939.2334
$ws.Range("I31").Value = 2539.0588
$ws.Range("J31").Value = 5206.8887
$ws.Range("K31").Value = 2539.0588
$ws.Range("L31").Value = 5206.8887
$ws.Range("M31").Value = -2244.0588
$ws.Range("N31").Value = -5796.8887

# Hunk 13: sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2939.2334
$ws.Range("I34").Value = 2539.0588
$ws.Range("J34").Value = 5206.8887
$ws.Range("K34").Value = 2539.0588
$ws.Range("L34").Value = 5206.8887
$ws.Range("M34").Value = -2337.0588
$ws.Range("N34").Value = -5610.8887

# Hunk 14: sheet CRP, row 64
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496

# Hunk 15: sheet CRP, row 67
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716

# Hunk 16: sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12906012
$ws.Range("I99").Value = 12906012
$ws.Range("K99").Value = 12906012
$ws.Range("M99").Value = -12904514

# Hunk 17: sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 31257680
$ws.Range("I107").Value = 52643784
$ws.Range("K107").Value = 52643784
$ws.Range("M107").Value = -52641864

# Hunk 18: sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 13412.8
$ws.Range("J122").Value = 3052.5
$ws.Range("L122").Value = 9157.5
$ws.Range("N122").Value = -14057.5

# Hunk 19: sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 12906012
$ws.Range("I126").Value = 12906012
$ws.Range("K126").Value = 38718036
$ws.Range("M126").Value = -38715566

# Hunk 20: sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5090.5
$ws.Range("I134").Value = 1386.6
$ws.Range("K134").Value = 4159.799999999999
$ws.Range("M134").Value = -1624.799999999999

# Hunk 21: sheet CRP, row 135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 74798.336
$ws.Range("J135").Value = 74798.336
$ws.Range("L135").Value = 74798.336
$ws.Range("N135").Value = -84938.336

# Hunk 22: sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 288208.1
$ws.Range("I5").Value = 2124.1667
$ws.Range("J5").Value = 347397.88
$ws.Range("K5").Value = 6372.500100000001
$ws.Range("L5").Value = 1042193.64
$ws.Range("M5").Value = -6260.500100000001
$ws.Range("N5").Value = -1042417.64

# Hunk 23: sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 288208.1
$ws.Range("I135").Value = 2124.1667
$ws.Range("J135").Value = 347397.88
$ws.Range("K135").Value = 19117.5003
$ws.Range("L135").Value = 3126580.92
$ws.Range("M135").Value = -16582.5003
$ws.Range("N135").Value = -3131650.92

# Hunk 24: sheet CUL, row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2420.5454
$ws.Range("I137").Value = 2265.8823
$ws.Range("J137").Value = 2946.4
$ws.Range("K137").Value = 6797.646900000001
$ws.Range("L137").Value = 8839.200000000001
$ws.Range("M137").Value = -1697.646900000001
$ws.Range("N137").Value = -19039.2

# Hunk 25: sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8202.27
$ws.Range("I102").Value = 10153.333
$ws.Range("K102").Value = 10153.333
$ws.Range("M102").Value = -8531.333000000001

# Hunk 26: sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 24213.783
$ws.Range("I7").Value = 54989.75
$ws.Range("K7").Value = 54989.75
$ws.Range("M7").Value = -54877.75

# Hunk 27: sheet LTW, row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5584.025
$ws.Range("I22").Value = 6898.136
$ws.Range("J22").Value = 3977.889
$ws.Range("K22").Value = 6898.136
$ws.Range("L22").Value = 3977.889
$ws.Range("M22").Value = -6603.136
$ws.Range("N22").Value = -4567.889

# Hunk 28: sheet LTW, row 24
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 30006
$ws.Range("I24").Value = 30006
$ws.Range("K24").Value = 30006
$ws.Range("M24").Value = -29663

# Hunk 29: sheet LTW, row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5584.025
$ws.Range("I27").Value = 6898.136
$ws.Range("J27").Value = 3977.889
$ws.Range("K27").Value = 6898.136
$ws.Range("L27").Value = 3977.889
$ws.Range("M27").Value = -6791.136
$ws.Range("N27").Value = -4191.889

# Hunk 30: sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 20899.965
$ws.Range("I40").Value = 23290.264
$ws.Range("K40").Value = 23290.264
$ws.Range("M40").Value = -23154.264

# Hunk 31: sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7258.243
$ws.Range("I122").Value = 5329.9062
$ws.Range("K122").Value = 15989.7186
$ws.Range("M122").Value = -13539.7186

# Hunk 32: sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 24213.783
$ws.Range("I126").Value = 54989.75
$ws.Range("K126").Value = 164969.25
$ws.Range("M126").Value = -162499.25

# Hunk 33: sheet LTW, row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 63332.332
$ws.Range("J140").Value = 63332.332
$ws.Range("L140").Value = 63332.332
$ws.Range("N140").Value = -73692.33199999999

# Hunk 34: sheet WVR, row 51
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 28097.8
$ws.Range("I51").Value = 6833
$ws.Range("K51").Value = 6833
$ws.Range("M51").Value = -6323

# Hunk 35: sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 15164.833
$ws.Range("I81").Value = 18219.777
$ws.Range("J81").Value = 6000
$ws.Range("K81").Value = 36439.554
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = -35378.554
$ws.Range("N81").Value = -14122

# Hunk 36: sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 15164.833
$ws.Range("I84").Value = 18219.777
$ws.Range("J84").Value = 6000
$ws.Range("K84").Value = 182197.77
$ws.Range("L84").Value = 60000
$ws.Range("M84").Value = -176893.77
$ws.Range("N84").Value = -70608

# Hunk 37: sheet WVR, row 98
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 29357.334

# Hunk 38: sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 56629
$ws.Range("I126").Value = 83608.39999999999
$ws.Range("J126").Value = 11663.333
$ws.Range("K126").Value = 250825.2
$ws.Range("L126").Value = 34989.999
$ws.Range("M126").Value = -248355.2
$ws.Range("N126").Value = -39929.999

# Hunk 39: sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13275.042
$ws.Range("I132").Value = 14100.068
$ws.Range("K132").Value = 42300.204
$ws.Range("M132").Value = -39770.204

# Hunk 40: sheet WVR, row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 145764.8
$ws.Range("J135").Value = 145764.8
$ws.Range("L135").Value = 145764.8
$ws.Range("N135").Value = -155904.8

# Hunk 41: sheet WVR, row 139
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 42825
$ws.Range("I139").Value = 42825
$ws.Range("K139").Value = 42825
$ws.Range("M139").Value = -37685

# Hunk 42: sheet WVR, row 141
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 99493
$ws.Range("J141").Value = 99493
$ws.Range("L141").Value = 99493
$ws.Range("N141").Value = -109853
